# Replaced MessagePack with simple byte structure for BLE packets
#
# On slide 3 (the "response APDU" BLE-packet figure):
#   - remove the "81" rectangle (shape "Rectangle 4") that used to label the
#     MessagePack message-type byte
#   - add three "2 bytes" labels above the remaining rectangles, matching the
#     style already used for the labels on slide 2

$p = $ppt.ActivePresentation
$s2 = $p.Slides.Item(2)
$s3 = $p.Slides.Item(3)

# EMU -> points, nudged by half an EMU so the COM layer's single-precision
# float round-trip lands back on the exact target EMU value instead of
# truncating it down by one.
function EmuToPt([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

# Drop the "81" rectangle (MessagePack message-type byte label).
$s3.Shapes.Item("Rectangle 4").Delete()

# The deck's shape-id counter for this slide is driven by how many shapes
# have been created on it so far (deleted or not). Burn through the three
# ids that would otherwise land on our new textboxes so the three "2 bytes"
# labels come out as ids 8, 10, 11 (names "TextBox 7"/"TextBox 9"/"TextBox 10"),
# matching the target document.
for ($i = 0; $i -lt 3; $i++) {
    $dummy = $s3.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $dummy.Delete()
}

# Reuse the existing "2 bytes" textbox from slide 2 as a style template
# (Consolas 10pt, no-wrap, auto-fit, no fill) via copy/paste, then move each
# copy into place and rename/retext it.
$template = $s2.Shapes.Item("TextBox 1")

$labels = @(
    @{ X = 1006615; Y = 1758310; Name = "TextBox 7" },
    @{ X = 1622576; Y = 1758310; Name = "TextBox 9" },
    @{ X = 2384967; Y = 1758310; Name = "TextBox 10" }
)

foreach ($label in $labels) {
    $template.Copy()
    $pasted = $s3.Shapes.Paste()
    $pasted.Name = $label.Name
    $pasted.Left = EmuToPt $label.X
    $pasted.Top = EmuToPt $label.Y
    $pasted.TextFrame.TextRange.Text = "2 bytes"
}
